# Commit: "Updated period to 2022Oct from 2021Oct"
#
# The "Map" sheet carries a period-tag column ("pe", column I) that was
# stamped with the prior reporting period "2021Oct" for every data row.
# Roll it forward to "2022Oct".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Map")
$ws.Activate()

$periodCol = 9  # column I ("pe")
$lastRow = $ws.Cells.Item($ws.Rows.Count, $periodCol).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $periodCol)
    if ($cell.Value() -eq "2021Oct") {
        $cell.Value = "2022Oct"
    }
}

# Best-effort view-state touch-up to mirror the author re-viewing the
# refreshed period column (zoom + selection), matching the saved workbook.
$excel.ActiveWindow.Zoom = 120
$ws.Range("I1:I1048576").Select() | Out-Null
